$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New status row for 29/1/2021, appended after row 21
$ws.Cells.Item(22, 1).Value = "29/1/2021"
$ws.Cells.Item(22, 2).Value = "1. Attended time management session`n2. Completed 2 bitwise operators programs, 1 string program and  2 recursion programs from the given list`n3. Worked on task given by Srinivas regarding testing on my mobile `n4. Completed one hacker rank program in C"
$ws.Cells.Item(22, 4).Value = "MaximizingXor.txt"

# Copy styles from row 21 so the new row matches formatting
$ws.Range("A21").Copy()
$ws.Range("A22").PasteSpecial(-4122)
$ws.Range("B21").Copy()
$ws.Range("B22").PasteSpecial(-4122)
$ws.Range("D21").Copy()
$ws.Range("D22").PasteSpecial(-4122)

$ws.Rows.Item(22).RowHeight = 105

$ws.Range("A22").Select()
